$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in columns A and B, rows 1-4
$ws.Range("A1").Value = -0.038633801308740008
$ws.Range("B1").Value = -0.041058936046815687

$ws.Range("A2").Value = -0.0068114147833885849
$ws.Range("B2").Value = -0.0070921184769611872

$ws.Range("A3").Value = -0.0034107977389385871
$ws.Range("B3").Value = -0.004190550261800958

$ws.Range("A4").Value = -0.013183342990520394
$ws.Range("B4").Value = -0.015146136896910507

# Add new row 5
$ws.Range("A5").Value = -0.065644490333685176
$ws.Range("B5").Value = -0.065633118481325123

# Adjust column widths for A and B (closest reachable value to the
# target stored width of 15.42578125 given this host's pixel-rounding
# of the ColumnWidth property)
$ws.Range("A1:B1").ColumnWidth = 14.67
